$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shp = $s.NotesPage.Shapes.Item(2)
$tf = $shp.TextFrame
$tr = $tf.TextRange
$t = $tr.Text
$idx = $t.IndexOf("laguna")
if ($idx -ge 0) {
    $sub = $tr.Characters($idx + 1, 6)
    $sub.Text = "alguna"
}
